# Companies!company "thecomp" is being rebranded to "Game On Wrestling Federation" (GOWF),
# a gameshow-themed wrestling promotion. Update the Companies, Bios and Notes sheets so every
# cell that shares the old strings picks up the new text (the engine dedupes the shared-string
# table by content on save).

$wb = $excel.ActiveWorkbook

$wsCompanies = $wb.Worksheets.Item("Companies")
$wsBios      = $wb.Worksheets.Item("Bios")
$wsNotes     = $wb.Worksheets.Item("Notes")

# ---- new field values -------------------------------------------------
$newUid             = 109
$newName            = '"Game On Wrestling Federation"'
$newInitials        = 'OWF'
$newUrl             = 'www."gameonwrestlingfederation".com'
$newLogo            = '"gameonwrestlingfederation".jpg'
$newBackdrop        = '"gameonwrestlingfederation"BD.jpg'
$newBanner          = '"gameonwrestlingfederation"Ban'
$newLogoDescShort   = 'a gameshow themed wrestling company'
$newLogoDescription = 'The logo for Game On Wrestling Federation would feature bold, vibrant colors, retro arcade fonts, and a mix of wrestling and game show symbols such as wrestling ring ropes intertwined with game controllers, dice, and a spotlight shining down on a wrestling championship belt.'

$newBio = @"
Name: Game On Wrestling Federation
Location: TBD
Size: Medium
Description:
Game On Wrestling Federation is a unique and innovative professional wrestling company that combines the thrill and excitement of a gameshow with the physicality and drama of professional wrestling. With a medium-sized roster of talented and charismatic wrestlers, GOWF offers a fresh and exciting take on sports entertainment.
Theme:
The theme of GOWF is all about competition and entertainment. Each wrestling event is designed like a gameshow, with different segments and matches that challenge the competitors in various ways. From trivia rounds to physical challenges, the wrestlers must showcase their skills not only in the ring but also in a range of different game formats.
Roster:
The roster of GOWF is diverse and eclectic, featuring a mix of seasoned veterans and up-and-coming talent. Wrestlers in the company are known for their athleticism, charisma, and willingness to take risks in order to entertain the audience. The roster includes a mix of traditional wrestlers, high-flying acrobats, and comedic characters, all of whom bring something unique to the ring.
Championships:
GOWF features a range of championships that are hotly contested by the wrestlers on the roster. Titles include the Game On Championship, which is the top prize in the company and is defended in high-stakes matches that combine in-ring action with gameshow-style challenges. Other titles include the Tag Team Championship, the Women's Championship, and the Hardcore Championship.
Events:
GOWF hosts regular wrestling events that are broadcast on television and live-streamed online for fans around the world to enjoy. Each event features a mix of traditional wrestling matches and gameshow-style challenges, with surprises and twists that keep the audience on the edge of their seats. From one-on-one grudge matches to multi-man spectacles, GOWF events are always full of excitement and drama.
Mission:
The mission of GOWF is to provide fans with a fresh and engaging alternative to traditional professional wrestling. By combining the excitement of a gameshow with the physicality of wrestling, GOWF offers a unique entertainment experience that appeals to a wide range of viewers. The company is dedicated to pushing the boundaries of what is possible in the world of sports entertainment and creating memorable moments that will keep fans coming back for more.
"@

# ---- Companies sheet (row 2 = the single company record) --------------
$wsCompanies.Range("A2").Value = $newUid
$wsCompanies.Range("B2").Value = $newName
$wsCompanies.Range("C2").Value = $newInitials
$wsCompanies.Range("D2").Value = $newUrl
$wsCompanies.Range("I2").Value = $newLogo
$wsCompanies.Range("J2").Value = $newBackdrop
$wsCompanies.Range("K2").Value = $newBanner
$wsCompanies.Range("M2").Value = 43
$wsCompanies.Range("R2").Value = 45

# ---- Bios sheet ---------------------------------------------------------
$wsBios.Range("A2").Value = $newUid
$wsBios.Range("B2").Value = $newBio
# Assigning a multi-line value auto-expands the row height (customHeight); put the
# row back the way it started (no explicit height) now that the new text is in place.
$wsBios.Rows.Item(2).AutoFit()

# ---- Notes sheet (shares Name / Logo / Backdrop / Banner / logo_description strings) ----
$wsNotes.Range("A2").Value = $newName
$wsNotes.Range("B2").Value = $newLogoDescShort
$wsNotes.Range("D2").Value = $newLogo
$wsNotes.Range("E2").Value = $newBackdrop
$wsNotes.Range("F2").Value = $newBanner
$wsNotes.Range("H2").Value = $newLogoDescription

Write-Host "Applied GOWF rebrand edits"
